$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix invalid facility utilisation values in rows 2-6 (columns E-H)
# Row 2
$ws.Range("E2").Value = 800000
$ws.Range("F2").Value = 761579.37

# Row 3
$ws.Range("E3").Value = 800000
$ws.Range("F3").Value = 761579.37

# Row 4
$ws.Range("E4").Value = 800000
$ws.Range("F4").Value = 761579.37

# Row 5
$ws.Range("E5").Value = 800000
$ws.Range("F5").Value = 761579.37
$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

# Row 6
$ws.Range("E6").Value = 800000
$ws.Range("G6").Value = 761579.37

# Align column G's width with columns E:F so E:G share one width/bestFit definition
$ws.Columns("G:G").ColumnWidth = 15.498697916666666

# Move / update the sheet selection to E2:H6 (active cell E2)
$excel.Goto($ws.Range("E2:H6"))
